# Update gh-pages to output generated at 456a3b4
# Refresh the view-count / favorite-count style numbers scraped for each
# 漫展 (convention) entry across the four sheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 191
$ws1.Range("F5").Value = 368
$ws1.Range("F8").Value = 816
$ws1.Range("F9").Value = 4214
$ws1.Range("F10").Value = 4214
$ws1.Range("F14").Value = 6104
$ws1.Range("F15").Value = 64
$ws1.Range("F17").Value = 2342
$ws1.Range("F21").Value = 9192
$ws1.Range("F22").Value = 43
$ws1.Range("F23").Value = 2473
$ws1.Range("F25").Value = 2314
$ws1.Range("F26").Value = 2456
$ws1.Range("F27").Value = 1394
$ws1.Range("F29").Value = 1970
$ws1.Range("F34").Value = 42
$ws1.Range("F37").Value = 58
$ws1.Range("F39").Value = 1221
$ws1.Range("F40").Value = 1218
$ws1.Range("F42").Value = 98
$ws1.Range("F43").Value = 241
$ws1.Range("F44").Value = 1541
$ws1.Range("F45").Value = 2533
$ws1.Range("F47").Value = 300
$ws1.Range("F48").Value = 1252

# ---------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 18
$ws2.Range("F5").Value = 169
$ws2.Range("G5").Value = 380
$ws2.Range("F22").Value = 72

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 896

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 191
$ws4.Range("F4").Value = 896
$ws4.Range("F6").Value = 368
$ws4.Range("F10").Value = 18
$ws4.Range("F11").Value = 169
$ws4.Range("G11").Value = 380
$ws4.Range("F13").Value = 816
$ws4.Range("F14").Value = 4214
$ws4.Range("F18").Value = 6104
$ws4.Range("F19").Value = 64
$ws4.Range("F23").Value = 9193
$ws4.Range("F24").Value = 43
$ws4.Range("F25").Value = 2473
$ws4.Range("F27").Value = 2314
$ws4.Range("F28").Value = 2456
$ws4.Range("F29").Value = 1394
$ws4.Range("F31").Value = 1970
$ws4.Range("F36").Value = 58
$ws4.Range("F38").Value = 1218
$ws4.Range("F40").Value = 98
$ws4.Range("F41").Value = 241
$ws4.Range("F42").Value = 1541
$ws4.Range("F43").Value = 2533
$ws4.Range("F45").Value = 300
$ws4.Range("F48").Value = 1252
$ws4.Range("F50").Value = 72
$ws4.Range("F51").Value = 72
